$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Ojaswi Sinha"
$ws.Range("A9").Value = "Yanye Luther"
$ws.Range("A8").Value = "James Yost"
$ws.Range("A11").Value = "Annie Rudnick"

$ws.Range("A12").Select()
